$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content so the sheet is rebuilt from scratch
$ws.Cells.Clear()

$ws.Cells.Item(1,1).Value = 'Method'
$ws.Cells.Item(1,2).Value = 'Endpoint'
$ws.Cells.Item(1,3).Value = 'Payload'

$ws.Cells.Item(2,1).Value = 'GET'
$ws.Cells.Item(2,2).Value = '/redfish/v1/AccountService'

$ws.Cells.Item(3,1).Value = 'GET'
$ws.Cells.Item(3,2).Value = '/redfish/v1/AccountService/Accounts/1'

$ws.Cells.Item(4,1).Value = 'GET'
$ws.Cells.Item(4,2).Value = '/redfish/v1/AccountService/Accounts/admin'

$ws.Cells.Item(5,1).Value = 'DELETE'
$ws.Cells.Item(5,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'

$ws.Cells.Item(6,1).Value = 'DELETE'
$ws.Cells.Item(6,2).Value = '/redfish/v1/AccountService/Accounts/${test2user.id}'

$ws.Cells.Item(7,1).Value = 'DELETE'
$ws.Cells.Item(7,2).Value = '/redfish/v1/AccountService/Accounts/${test3user.id}'

$ws.Cells.Item(8,1).Value = 'POST'
$ws.Cells.Item(8,2).Value = '/redfish/v1/AccountService/Accounts'
$ws.Cells.Item(8,3).Value = '{"UserName": "testuser", "Password": "Test1234", "RoleId": "Administrator"}'

$ws.Cells.Item(9,1).Value = 'POST'
$ws.Cells.Item(9,2).Value = '/redfish/v1/AccountService/Accounts'
$ws.Cells.Item(9,3).Value = '{"UserName": "testuser", "Password": "Test1234", "RoleId": "Administrator"}'

$ws.Cells.Item(10,1).Value = 'POST'
$ws.Cells.Item(10,2).Value = '/redfish/v1/AccountService/Accounts'
$ws.Cells.Item(10,3).Value = '{"UserName": "test2user", "Password": "Test1234", "RoleId": "Administrator"}'

$ws.Cells.Item(11,1).Value = 'POST'
$ws.Cells.Item(11,2).Value = '/redfish/v1/AccountService/Accounts'
$ws.Cells.Item(11,3).Value = '{"UserName": "test3user", "Password": "Test1234", "RoleId": "Administrator","redundant-key":"redundant-value"}'

$ws.Cells.Item(12,1).Value = 'PATCH'
$ws.Cells.Item(12,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(12,3).Value = '{"Password":"simplepasswordtest"}'

$ws.Cells.Item(13,1).Value = 'PATCH'
$ws.Cells.Item(13,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(13,3).Value = '{"Password":"aD0!"}'

$ws.Cells.Item(14,1).Value = 'PATCH'
$ws.Cells.Item(14,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(14,3).Value = '{"Password":"Vveryvemmmmmmmmmmmmmmmmmmmmmmmmmryjjjjjjjjjjjjjjjjjjjjveryveryveryveryveryveryveryveryveryveryveryveryveryveryveryveryveryverylongpassword012!"}'

$ws.Cells.Item(15,1).Value = 'PATCH'
$ws.Cells.Item(15,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(15,3).Value = '{"Password":"Test5678!","RoleId":"Manager"}'

$ws.Cells.Item(16,1).Value = 'PATCH'
$ws.Cells.Item(16,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(16,3).Value = '{"Password":"Test5678!"}'

$ws.Cells.Item(17,1).Value = 'PATCH'
$ws.Cells.Item(17,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(17,3).Value = '{"Password":"Test5678!","RoleId":"Manager","Testkey1":"testvalue"}'

$ws.Cells.Item(18,1).Value = 'POST'
$ws.Cells.Item(18,2).Value = '/redfish/v1/SessionService/Sessions'
$ws.Cells.Item(18,3).Value = '{"UserName":"test2user","Password":"Test1234"}'

$ws.Cells.Item(19,1).Value = 'PATCH'
$ws.Cells.Item(19,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(19,3).Value = '{"Password":"Testxxxx!","TestKey1":"TestValue1"}'

$ws.Cells.Item(20,1).Value = 'PATCH'
$ws.Cells.Item(20,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(20,3).Value = '{"Password":"Test9012!","TestKey1":"TestValue1","RoleId":"Operator"}'

$ws.Cells.Item(21,1).Value = 'PATCH'
$ws.Cells.Item(21,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(21,3).Value = '{"Password":"TestAABB!","TestKey3":"TestValue3","TestKey4":"TestValue4"}'

$ws.Cells.Item(22,1).Value = 'PATCH'
$ws.Cells.Item(22,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'
$ws.Cells.Item(22,3).Value = '{"TestKey5":"TestValue5","TestKey6":"TestValue6"}'

$ws.Cells.Item(23,1).Value = 'DELETE'
$ws.Cells.Item(23,2).Value = '/redfish/v1/AccountService/Accounts/${testuser.id}'

$ws.Cells.Item(24,1).Value = 'DELETE'
$ws.Cells.Item(24,2).Value = '/redfish/v1/AccountService/Accounts/${test2user.id}'

$ws.Cells.Item(25,1).Value = 'DELETE'
$ws.Cells.Item(25,2).Value = '/redfish/v1/AccountService/Accounts/${test3user.id}'

$ws.Cells.Item(26,1).Value = 'DELETE'
$ws.Cells.Item(26,2).Value = '/redfish/v1/AccountService/Accounts/testwronguser'

$ws.Cells.Item(27,1).Value = 'DELAY'
$ws.Cells.Item(27,2).Value = 5

# --- Formatting to match the target layout ---

# Column C is much wider now to comfortably fit long JSON payloads
$ws.Columns.Item(3).ColumnWidth = 143.5

# Row 18 (the new SessionService POST) is taller and its payload cell wraps text
$ws.Rows.Item(18).RowHeight = 16
$ws.Range("C18").WrapText = $true

# Restore the selection to match the saved workbook state
$ws.Range("C22").Select() | Out-Null
